$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194, shifting existing rows 194-233 down to 195-234.
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new record.
$ws.Range("A194").Value = 4
$ws.Range("B194").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C194").Value = "Los Lagos"
$ws.Range("D194").Value = 44637
$ws.Range("E194").Value = 10
$ws.Range("F194").Value = 100112003
$ws.Range("G194").Value = "Ajo"
$ws.Range("H194").Value = "Chino"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 80
$ws.Range("K194").Value = 21000
$ws.Range("L194").Value = 22000
$ws.Range("M194").Value = 21500
$ws.Range("N194").Value = "$/caja 10 kilos"
$ws.Range("O194").Value = "China"
$ws.Range("P194").Value = 2150
$ws.Range("Q194").Value = 10
$ws.Range("R194").Value = "Hortaliza"
